$wb = $excel.ActiveWorkbook

# --- BDSBaPCF sheet: update boolean flags ---
$bds = $wb.Worksheets.Item("BDSBaPCF")

# hard coal: 1 -> 0
$bds.Range("B2").Value = 0
# nuclear: 1 -> 0
$bds.Range("B4").Value = 0
# biomass: 0 -> 1
$bds.Range("B9").Value = 1
# (B13 = lignite =B2, B17 = municipal solid waste =B9 recalc automatically)

# --- About sheet: append explanatory notes about coal in the US ---
$about = $wb.Worksheets.Item("About")

$about.Range("A24").Value = "For the United States, we have set coal to 0 as of version 3.4. This reflects"
$about.Range("A25").Value = "the fact that certain air quality / environmental restrictions, as well as current"
$about.Range("A26").Value = "supply chain logistics, limit the amount the coal dispatches annually. "

# --- Update selections to match the saved workbook state ---
$bds.Range("B4").Select()
$about.Range("A27").Select()
